$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 0.8
$ws.Range("B4").Value = 0.7

$ws.Range("C13").Select() | Out-Null
